# Append two new donation records (rows 13 and 14) to the LINS-Donations sheet,
# mirroring the existing data layout (text-typed receipt/contact fields, numeric
# row-number column). Columns whose content is all-digits (F, H, L) are entered
# with a leading apostrophe so Excel stores them as text, matching every other
# "number looking" value already in the sheet (Donation Amount, Phone, Zip).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 ---------------------------------------------------------------
$ws.Cells.Item(13, 1).Value  = "REC-1741892275029-290"
$ws.Cells.Item(13, 2).Value  = 12
$ws.Cells.Item(13, 3).Value  = "2025-03-13T18:57:55.037Z"
$ws.Cells.Item(13, 4).Value  = "Deepak"
$ws.Cells.Item(13, 5).Value  = "Adhikari"
$ws.Cells.Item(13, 6).Value  = "'123456"
$ws.Cells.Item(13, 7).Value  = "longislandnepalese@gmail.com"
$ws.Cells.Item(13, 8).Value  = "'3477712375"
$ws.Cells.Item(13, 9).Value  = "11 alpine ln"
$ws.Cells.Item(13, 10).Value = "Hicksville"
$ws.Cells.Item(13, 11).Value = "NY"
$ws.Cells.Item(13, 12).Value = "'11801"

# --- Row 14 ---------------------------------------------------------------
$ws.Cells.Item(14, 1).Value  = "REC-1741892978780-739"
$ws.Cells.Item(14, 2).Value  = 13
$ws.Cells.Item(14, 3).Value  = "2025-03-13T19:09:38.787Z"
$ws.Cells.Item(14, 4).Value  = "Deepak"
$ws.Cells.Item(14, 5).Value  = "Adhikari"
$ws.Cells.Item(14, 6).Value  = "'2222"
$ws.Cells.Item(14, 7).Value  = "dadhikari856@gmail.com"
$ws.Cells.Item(14, 8).Value  = "'3477712375"
$ws.Cells.Item(14, 9).Value  = "11 alpine ln"
$ws.Cells.Item(14, 10).Value = "Hicksville"
$ws.Cells.Item(14, 11).Value = "NY"
$ws.Cells.Item(14, 12).Value = "'11801"
